$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "2025/10/08"
$ws.Range("A76").ClearFormats()
$ws.Range("B76").Value = "水"
$ws.Range("C76").Value = 2
$ws.Range("D76").Value = 14
